$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update cell values (Package / Url / Notes table)
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "Package"
$ws.Range("B1").Value = "Url"
$ws.Range("C1").Value = "Notes"

$ws.Range("A2").Value = "lfe"
$ws.Range("B2").Value = "https://cran.r-project.org/package=lfe"
$ws.Range("C2").Value = "This package is a workhorse package for regression modeling. One of its chief strengths is that it allows for the  fast estimation of models with high dimensional fixed effects. It also has functionality for instrumental variables regression and allows for heteroskedastic/cluster robust standard errors.  "

$ws.Range("A3").Value = "estimatr"
$ws.Range("B3").Value = "https://cran.r-project.org/package=estimatr"
$ws.Range("C3").Value = "This package has my goto drop-in replacement for ``lm``: ``robust_lm`` which conveniently allows for robust and clustered standard errors. Also includes other estimators commonly used in designed-based inference. "

$ws.Range("A4").Value = "sensemakr"
$ws.Range("B4").Value = "https://cran.r-project.org/package=sensemakr"
$ws.Range("C4").Value = "R package implementing the sensitivity analysis methods for unmeasured confounding proposed in [Cinelli and Hazlett](https://doi.org/10.1111/rssb.12348). This is the first approach I suggest to students interested in implementing a sensitivity analysis. "

$ws.Range("A5").Value = "rdrobust"
$ws.Range("B5").Value = "https://cran.r-project.org/package=rdrobust"
$ws.Range("C5").Value = "Implements a suite of estimation methods, bandwidth selection algorthims, and graphical tools for regression discontinuity designs. "

$ws.Range("A6").Value = "gsynth"
$ws.Range("B6").Value = "https://cran.r-project.org/package=gsynth"
$ws.Range("C6").Value = "Provides recent panel-data estimators that go beyond standard fixed effects models, including interactive fixed effects models and matrix completion methods. "

$ws.Range("A7").Value = "PanelMatch"
$ws.Range("B7").Value = "https://cran.r-project.org/package=PanelMatch"
$ws.Range("C7").Value = "Implements generalized diffence-in-differences  estimators that avoid some of the recently identified [problems](https://www.nber.org/papers/w25018) with two-way fixed effect models. "

# ---------------------------------------------------------------------------
# 2. Column widths / formats
#    Column B holds URLs -> store as text (numFmt "@"), default width.
#    Column C holds long notes -> widen considerably.
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).NumberFormat = "@"
$ws.Columns.Item(3).ColumnWidth = 84.92

# ---------------------------------------------------------------------------
# 3. Hyperlinks - drop the stale one and recreate the full set that lines up
#    with the new rows.
# ---------------------------------------------------------------------------

$ws.Range("B3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://cran.r-project.org/package=lfe") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://cran.r-project.org/package=estimatr") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://cran.r-project.org/package=rdrobust") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://cran.r-project.org/package=gsynth") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://cran.r-project.org/package=PanelMatch") | Out-Null

# ---------------------------------------------------------------------------
# 4. Sheet view / selection / page setup touch-ups
# ---------------------------------------------------------------------------

$ws.Range("C8").Select() | Out-Null
$ws.PageSetup.Orientation = 1

Write-Output "applied useful_r_packages update"
